# Commit: "bug fix - valid space"
# Adds a new "volume" column (J) to Sheet1, computed as the bounding-box
# volume of each part: (X-Max - X-Min) * (Y-Max - Y-Min) * (Z-Max - Z-Min),
# using the already-present X/Y/Z-Min/Max_transf columns (D..I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell J1: "volume", styled like the other header cells
#     (bold font, thin border, centered horizontally, top-aligned). ---
$ws.Range("J1").Value = "volume"
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").Borders.LineStyle = 1
$ws.Range("J1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("J1").VerticalAlignment = -4160     # xlTop

# --- Data rows 2..165: bounding-box volume for each of the 164 parts. ---
$volumes = @(
    3116867419.8635159,12286387.26057039,287956071.11179131,133447286.8565042,377079640.54395348,4395391.1287168721,460341926.85456371,3005053.9451211751,
    3005053.9451211751,29217989.5041558,7778674.0579851121,70853517.605508074,270554250.13993752,24197285.633770879,190151379.48631811,12286387.26057039,
    207108548.56353459,333926592.1432296,172579455.72286719,377079640.54395348,2861444.451779569,460341926.85456371,262488333.85280579,2605736.7509720302,
    31181455.210011311,533956456.56468308,2917773.1389583428,7760019.7748016911,75628347.865811318,188608651.60712871,24197285.633770879,6662053.7555670291,
    3308752118.5323648,21358983.239237592,304351430.39889699,141943383.915791,106813908.8100172,1663192.237470553,4103728.501684946,31133912.31276805,
    6234212.0570974071,28665715.39918004,158155300.7097275,23444983.158251811,4929621.1115717953,401903963.47719991,21151233.43516472,222806223.4113988,
    224542997.25596669,303747215.79968822,92183437.884501725,575403670.1906333,46058035.791716702,995518.69829942635,573162146.54022479,180270454.24290681,
    6776060.7104884097,4488384.5768924421,1392925.6835431389,1735465.935423509,35329731.61328961,771034612.65438604,22332724.49274718,18854542.713296209,
    60345431.267341457,177435737.996445,9293227.648824852,1655607.6914827011,121156543.33588611,134824515.9012163,8041068.4114933386,1788924.163083768,
    502801.10472188448,39937537.631623253,371403503.76864201,22362241.181072552,235212960.05552709,293994355.63744599,39858456.13622795,663872194.80202699,
    484270405.14201868,1005639.710654537,646682464.48011625,184541291.3104502,8558267.4259310309,1893027.3603401461,1773577.166692896,35306845.062677547,
    750528201.86815631,12548810.728767401,40388034.09709233,21163565.08271651,528223735.99014068,188246074.5602878,19960114.04006717,5491812.4317981647,
    130158596.4370307,121300707.7582356,7385317.1001000712,112425202.7303035,263156444.7873379,12845468.493239909,223885729.01405829,212292401.0615741,
    278851087.69414783,436881020.73227191,75781424.474848121,544690634.14361703,120543018.34000351,995504.56402691628,743385241.11292827,151156203.30609179,
    12987533.297290999,4383993.8856919464,46875677.065853037,764153626.45658159,51057387.043277957,72584427.855768412,80908233.498167738,156614795.32033539,
    19763764.451679971,1536995.6537733199,113623993.7086636,115667993.73727059,6905294.4530950673,12845468.493239909,200992712.0210312,171542669.06311911,
    247125115.75061259,120543018.34000351,11343077.455338471,746775051.3485837,558414593.87999749,2140993.721058113,724877236.05331504,322736281.40574932,
    12987533.297290999,51057387.043277957,74089778.389638901,83927342.146729976,156249925.68286541,19084242.942981709,4539732.320076854,46744007.495347701,
    772801582.2698251,1536995.6537733199,87536796.362904146,101185162.6066363,399961538.39839971,6905294.4530950673,0,0,
    0,0,0,0,0,3007265.5118336622,2709722.554991364,512050.74384261732,
    748721.45857866819,1757254.837009503,1742094.473237742,1618661.12250096
)

$firstRow = 2
$lastRow = 165
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 10).Value = $volumes[$r - $firstRow]
}

# --- Restore the selection state recorded in the saved workbook. ---
$ws.Range("A2:XFD4087").Select()
